$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.41829962205201809
$ws.Range("B1").Value = 0.41725173969997797
$ws.Range("A2").Value = -0.29251371723288777
$ws.Range("B2").Value = 0.28944809341118827
$ws.Range("A3").Value = -0.18649682980883142
$ws.Range("B3").Value = 0.18558314544700494
$ws.Range("A4").Value = -0.17358314564053501
$ws.Range("B4").Value = 0.17276195128539484
$ws.Range("A5").Value = -0.16676195197336163
$ws.Range("B5").Value = 0.16511814669234504
$ws.Range("A6").Value = -0.063796773441434862
$ws.Range("B6").Value = 0.063734715822217236
$ws.Range("A7").Value = -0.043734716656027572
$ws.Range("B7").Value = 0.043626282188958143
$ws.Range("A8").Value = -0.023626283027608608
$ws.Range("B8").Value = 0.023581436819195112
$ws.Range("A9").Value = -0.01758143754095709
$ws.Range("B9").Value = 0.017555890908484528
$ws.Range("A10").Value = -0.01155589163299453
$ws.Range("B10").Value = 0.011557360584973253
$ws.Range("A11").Value = -0.0070573612966420285
$ws.Range("B11").Value = 0.0070579161191730577
$ws.Range("A12").Value = -0.0010579168439033282
$ws.Range("B12").Value = 0.001057450051869413
$ws.Range("A13").Value = 0.0049425492232746393
$ws.Range("B13").Value = -0.004943794441584437
$ws.Range("A14").Value = -0.027084313177473618
$ws.Range("B14").Value = 0.027052462722002346
$ws.Range("A15").Value = -0.021052463448884673
$ws.Range("B15").Value = 0.021027430496753929
$ws.Range("A16").Value = -0.015027431225929089
$ws.Range("B16").Value = 0.015004133424103117
$ws.Range("A17").Value = -0.0090041341563242838
$ws.Range("B17").Value = 0.0089999992412268526
$ws.Range("A18").Value = -0.10467052511267738
$ws.Range("B18").Value = 0.10455818217332435
$ws.Range("A19").Value = -0.027096529700363359
$ws.Range("B19").Value = 0.027013238939074924
$ws.Range("A20").Value = -0.018013239618458243
$ws.Range("B20").Value = 0.018004241549915534
$ws.Range("A21").Value = -0.0090042422302305525
$ws.Range("B21").Value = 0.0089999993190659211
$ws.Range("A22").Value = -0.093951264296302739
$ws.Range("B22").Value = 0.093637366467813621
$ws.Range("A23").Value = -0.084637367167717414
$ws.Range("B23").Value = 0.084127281352007088
$ws.Range("A24").Value = -0.042127282347874662
$ws.Range("B24").Value = 0.04199999899877227
$ws.Range("A25").Value = -0.095321373962114819
$ws.Range("B25").Value = 0.095198334881924751
$ws.Range("A26").Value = -0.089198335595739309
$ws.Range("B26").Value = 0.089044727947559466
$ws.Range("A27").Value = -0.083044728664873002
$ws.Range("B27").Value = 0.082536368138486527
$ws.Range("A28").Value = -0.077416825368148423
$ws.Range("B28").Value = 0.077070990531693973
$ws.Range("A29").Value = -0.065070991325693939
$ws.Range("B29").Value = 0.064921008968472194
$ws.Range("A30").Value = -0.042169241747813668
$ws.Range("B30").Value = 0.042018995545197058
$ws.Range("A31").Value = -0.027018996377293547
$ws.Range("B31").Value = 0.027000603728595252
$ws.Range("A32").Value = -0.0060006046136749092
$ws.Range("B32").Value = 0.0059999992440955552
